# JobMaterial.xlsx — re-sequence the "Folded Sheet" element numbering.
#
# The sheet lists Ink/Varnish, Plate and Sheet material rows for three
# "Folded Sheet" elements (unsuffixed, " 1" and " 2"). The edit rotates
# which physical sheet each label refers to (unsuffixed -> 1 -> 2 -> back
# to unsuffixed), which also reshuffles which CMYK ink/plate/quantity
# values land on each row.
#
# Written as direct cell writes (rather than shared-string-table surgery)
# since that is the supported COM surface and is robust to how the host
# re-serialises the shared strings table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ink / Varnish rows (2-13): element label + ink colour/weight shuffle ---
$ws.Range("B2").Value = "Folded Sheet  4p"
$ws.Range("D2").Value = "Black - Sheet-fed Offset - "
$ws.Range("E2").Value = "0.07"
$ws.Range("H2").Value = "M594GK - Black (General)"

$ws.Range("B3").Value = "Folded Sheet  4p"
$ws.Range("D3").Value = "Yellow - Sheet-fed Offset - "
$ws.Range("E3").Value = "0.07"

$ws.Range("B4").Value = "Folded Sheet  4p"
$ws.Range("E4").Value = "0.07"

$ws.Range("B5").Value = "Folded Sheet  4p"
$ws.Range("D5").Value = "Magenta - Sheet-fed Offset - "
$ws.Range("E5").Value = "0.07"
$ws.Range("H5").Value = "M594GN - 4/C Process (General)"

$ws.Range("B6").Value = "Folded Sheet 1  4p"
$ws.Range("D6").Value = "Yellow - Sheet-fed Offset - "

$ws.Range("B7").Value = "Folded Sheet 1  4p"

$ws.Range("B8").Value = "Folded Sheet 1  4p"
$ws.Range("D8").Value = "Magenta - Sheet-fed Offset - "

$ws.Range("B9").Value = "Folded Sheet 1  4p"
$ws.Range("D9").Value = "Cyan - Sheet-fed Offset - "

$ws.Range("B10").Value = "Folded Sheet 2  4p"
$ws.Range("D10").Value = "Yellow - Sheet-fed Offset - "
$ws.Range("E10").Value = "0.08"

$ws.Range("B11").Value = "Folded Sheet 2  4p"
$ws.Range("E11").Value = "0.08"

$ws.Range("B12").Value = "Folded Sheet 2  4p"
$ws.Range("D12").Value = "Magenta - Sheet-fed Offset - "
$ws.Range("E12").Value = "0.08"

$ws.Range("B13").Value = "Folded Sheet 2  4p"
$ws.Range("E13").Value = "0.08"

# --- Plate rows (14-16): just the element suffix rotates ---
$ws.Range("B14").Value = "Plate - Folded Sheet 2  4p"
$ws.Range("B15").Value = "Plate - Folded Sheet  4p"
$ws.Range("B16").Value = "Plate - Folded Sheet 1  4p"

# --- Sheet (paper) rows (17-19): element suffix + quantity rotate together ---
$ws.Range("B17").Value = "Folded Sheet  4p"
$ws.Range("E17").Value = "1,001.00"

$ws.Range("B18").Value = "Folded Sheet 1  4p"
$ws.Range("E18").Value = "1,106.00"

$ws.Range("B19").Value = "Folded Sheet 2  4p"
$ws.Range("E19").Value = "1,213.00"
